$d = $word.ActiveDocument

# 1. Merge the title's two runs into a single new headline.
$null = $d.Content.Find.Execute("SEATTLE, KNOW YOUR WORTH; BEFORE YOU POST TO HOST", $true, $false, $false, $false, $false, $true, 1, $false, "BEFORE YOU POST TO HOST: KNOW YOUR WORTH", 2)

# 2. Rewrite the "AirBnB is a great platform..." sentence with a shorter one, and
#    split what used to follow it into a brand-new paragraph with new copy.
$oldIntro = " is a great platform to get your rental property out on the market, whether you are a full time vacation home renter, or just looking to rent out your loft while you are out of town.   That being said, there are a lot of properties postings competing for your business; post too high and you will price out some renters, post too low and you are missing an opportunity to earn your properties potential."
$newIntro = " is a great way to vehicle to help you rent your property`rBut with so many people out there, how I can get the most money for what I have."
$null = $d.Content.Find.Execute($oldIntro, $true, $false, $false, $false, $false, $true, 1, $false, $newIntro, 2)

# 3. Only one blank paragraph should remain between the new paragraph and "Let's ...".
#    There are currently three; delete two of them.
$d.Paragraphs(5).Range.Delete()
$d.Paragraphs(5).Range.Delete()

# 4. Split "Let's look to answer a few questions" after "Let's" and drop the
#    _GoBack bookmark there (it moves up from further down in the document).
$pLets = $d.Paragraphs(6)
$letsStart = $pLets.Range.Start
$bmRange = $d.Range($letsStart + 5, $letsStart + 5)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 5. Clear out the first bulleted list item's text (keep the paragraph/pPr so the
#    list-numbering paragraph stays behind, now empty).
$pBullet = $d.Paragraphs(7)
$bulletTextRange = $d.Range($pBullet.Range.Start, $pBullet.Range.End - 1)
$bulletTextRange.Delete()

# 6. Remove every paragraph after that now-empty bullet item through the end of
#    the body (the other two bullets, the body copy, and all embedded images),
#    leaving the empty bullet paragraph as the document's last paragraph.
while ($d.Paragraphs.Count -gt 7) {
    $d.Paragraphs(8).Range.Delete()
}

# 7. Drop the now-unused "Balloon Text" style pair from the style sheet (delete
#    the character style before the paragraph style it is linked to).
$d.Styles("Balloon Text Char").Delete()
$d.Styles("Balloon Text").Delete()
